$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L5").Value = 1.6
$ws.Range("N5").Value = 3.1
$ws.Range("O4").Value = 2021
$ws.Range("O5").Value = 4.1

$ws.Range("P4").Select()
